$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reorder element/metal column labels M1:V1 ---
$ws.Range("M1").Value2 = "Cu"
$ws.Range("N1").Value2 = "Mn"
$ws.Range("O1").Value2 = "As"
$ws.Range("P1").Value2 = "Se"
$ws.Range("Q1").Value2 = "Ni"
$ws.Range("R1").Value2 = "Co"
$ws.Range("S1").Value2 = "Sr"
$ws.Range("T1").Value2 = "Cd"
$ws.Range("U1").Value2 = "Pb"
$ws.Range("V1").Value2 = "Ag"

# --- Data rows 2-35: permute trace-element columns M:V to match new header order ---
# row 2
$ws.Range("M2").Value2 = 0.0000644919674844881
$ws.Range("N2").Value2 = 0.0000434080550376362
$ws.Range("O2").Value2 = 0.0000302434102988301
$ws.Range("P2").Value2 = 0.000029651740872367
$ws.Range("Q2").Value2 = 0.0000144503879155434
$ws.Range("R2").Value2 = 0.00000125160840213368
$ws.Range("S2").Value2 = 0.00241321478189574
$ws.Range("T2").Value2 = 0.00000157019963176771
$ws.Range("U2").Value2 = 0.000000534778135457117
$ws.Range("V2").Value2 = 0.000000113782582012153
# row 3
$ws.Range("M3").Value2 = 0.0000787531644368318
$ws.Range("N3").Value2 = 0.0000507187838107039
$ws.Range("O3").Value2 = 0.000121428275870377
$ws.Range("P3").Value2 = 0.0000331349418699386
$ws.Range("Q3").Value2 = 0.000014391314609676
$ws.Range("R3").Value2 = 0.00000143414313701277
$ws.Range("S3").Value2 = 0.00208439610614425
$ws.Range("T3").Value2 = 0.0000109867835105065
$ws.Range("U3").Value2 = 0.000000785661022885258
$ws.Range("V3").Value2 = 0.000000399065916386163
# row 4
$ws.Range("M4").Value2 = 0.000136516253505943
$ws.Range("N4").Value2 = 0.0000397722670492563
$ws.Range("O4").Value2 = 0.0000917695200174876
$ws.Range("P4").Value2 = 0.0000471870377793832
$ws.Range("Q4").Value2 = 0.0000151345794966197
$ws.Range("R4").Value2 = 0.00000119668768112807
$ws.Range("S4").Value2 = 0.00140050001834921
$ws.Range("T4").Value2 = 0.0000194989698630867
$ws.Range("U4").Value2 = 0.00000307384482799562
$ws.Range("V4").Value2 = 0.00000175983482518833
# row 5
$ws.Range("M5").Value2 = 0.000134288476190443
$ws.Range("N5").Value2 = 0.0000372782004923443
$ws.Range("O5").Value2 = 0.0000210794339127139
$ws.Range("P5").Value2 = 0.0000260442506825747
$ws.Range("Q5").Value2 = 0.0000103083060052194
$ws.Range("R5").Value2 = 0.00000084149436777301
$ws.Range("S5").Value2 = 0.00146051164961297
$ws.Range("T5").Value2 = 0.00000524531489245176
$ws.Range("U5").Value2 = 0.000000434772090016055
$ws.Range("V5").Value2 = 0.000000182323779684152
# row 6
$ws.Range("M6").Value2 = 0.0000560967026090874
$ws.Range("N6").Value2 = 0.0000116418983124163
$ws.Range("O6").Value2 = 0.000136210210255271
$ws.Range("P6").Value2 = 0.0000549756309197436
$ws.Range("Q6").Value2 = 0.0000115844074565525
$ws.Range("R6").Value2 = 0.00000110669897537785
$ws.Range("S6").Value2 = 0.00107020665461831
$ws.Range("T6").Value2 = 0.0000260146122783623
$ws.Range("U6").Value2 = 0.000000201217995523245
$ws.Range("V6").Value2 = 0.000000273081565352975
# row 7
$ws.Range("M7").Value2 = 0.0000560522395121846
$ws.Range("N7").Value2 = 0.0000198637834058653
$ws.Range("O7").Value2 = 0.000172577109773352
$ws.Range("P7").Value2 = 0.0000486243036047801
$ws.Range("Q7").Value2 = 0.0000134150480888907
$ws.Range("R7").Value2 = 0.00000116105212865269
$ws.Range("S7").Value2 = 0.00096728232038646
$ws.Range("T7").Value2 = 0.0000034551792262315
$ws.Range("U7").Value2 = 0.00000015387437849614
$ws.Range("V7").Value2 = 0.000000125897218769569
# row 8
$ws.Range("M8").Value2 = 0.0000795626852444338
$ws.Range("N8").Value2 = 0.0000373492031211244
$ws.Range("O8").Value2 = 0.0000391680204871588
$ws.Range("P8").Value2 = 0.0000540992421199526
$ws.Range("Q8").Value2 = 0.0000174691063295859
$ws.Range("R8").Value2 = 0.00000148043273979542
$ws.Range("S8").Value2 = 0.00102333855686401
$ws.Range("T8").Value2 = 0.00000444129821938626
$ws.Range("U8").Value2 = 0.00000167077409205483
$ws.Range("V8").Value2 = 0.000000401831743658757
# row 9
$ws.Range("M9").Value2 = 0.00008320461609182
$ws.Range("N9").Value2 = 0.0000395919634274323
$ws.Range("O9").Value2 = 0.000125524341236873
$ws.Range("P9").Value2 = 0.0000437230245257944
$ws.Range("Q9").Value2 = 0.0000145217834793569
$ws.Range("R9").Value2 = 0.00000168711273864407
$ws.Range("S9").Value2 = 0.00160034468816969
$ws.Range("T9").Value2 = 0.00000633849832649453
$ws.Range("U9").Value2 = 0.000000788370438618723
$ws.Range("V9").Value2 = 0.00000143483419828608
# row 10
$ws.Range("M10").Value2 = 0.0000913628549463918
$ws.Range("N10").Value2 = 0.0000498755039409583
$ws.Range("O10").Value2 = 0.0000754681867712682
$ws.Range("P10").Value2 = 0.0000803549785715439
$ws.Range("Q10").Value2 = 0.0000155672027452082
$ws.Range("R10").Value2 = 0.00000136024101657159
$ws.Range("S10").Value2 = 0.00182773570224848
$ws.Range("T10").Value2 = 0.0000326709740461732
$ws.Range("U10").Value2 = 0.00000103277558665621
$ws.Range("V10").Value2 = 0.000000604551562920707
# row 11
$ws.Range("M11").Value2 = 0.000099417261660419
$ws.Range("N11").Value2 = 0.000055860689101275
$ws.Range("O11").Value2 = 0.0000891433243477396
$ws.Range("P11").Value2 = 0.0000497393911215648
$ws.Range("Q11").Value2 = 0.0000189483394748818
$ws.Range("R11").Value2 = 0.00000163029544183237
$ws.Range("S11").Value2 = 0.00116790674217908
$ws.Range("T11").Value2 = 0.0000394962141002407
$ws.Range("U11").Value2 = 0.00000209169981216228
$ws.Range("V11").Value2 = 0.00000239930272571556
# row 12
$ws.Range("M12").Value2 = 0.0000990493895828377
$ws.Range("N12").Value2 = 0.0000210824939127557
$ws.Range("O12").Value2 = 0.0000744863675841216
$ws.Range("P12").Value2 = 0.0000408003754922494
$ws.Range("Q12").Value2 = 0.0000182766056320035
$ws.Range("R12").Value2 = 0.00000173259768155738
$ws.Range("S12").Value2 = 0.00175236156130541
$ws.Range("T12").Value2 = 0.0000191352381113594
$ws.Range("U12").Value2 = 0.00000147194139318149
$ws.Range("V12").Value2 = 0.00000410916972263166
# row 13
$ws.Range("M13").Value2 = 0.0000838425905988663
$ws.Range("N13").Value2 = 0.0000352003204761797
$ws.Range("O13").Value2 = 0.0000374867081730622
$ws.Range("P13").Value2 = 0.0000503759047499931
$ws.Range("Q13").Value2 = 0.0000333410601512423
$ws.Range("R13").Value2 = 0.00000206026144114685
$ws.Range("S13").Value2 = 0.00119791640232926
$ws.Range("T13").Value2 = 0.0000105274334614699
$ws.Range("U13").Value2 = 0.00000429639885897697
$ws.Range("V13").Value2 = 0.000000376877092892716
# row 14
$ws.Range("M14").Value2 = 0.000074675815157623
$ws.Range("N14").Value2 = 0.0000445388618941192
$ws.Range("O14").Value2 = 0.0000511439448144716
$ws.Range("P14").Value2 = 0.0000814828883813501
$ws.Range("Q14").Value2 = 0.0000266425210151217
$ws.Range("R14").Value2 = 0.00000216129624610919
$ws.Range("S14").Value2 = 0.00117245271593839
$ws.Range("T14").Value2 = 0.000010766083169871
$ws.Range("U14").Value2 = 0.00000177751466969727
$ws.Range("V14").Value2 = 0.000000403980606749381
# row 15
$ws.Range("M15").Value2 = 0.0000885642628478717
$ws.Range("N15").Value2 = 0.000040985213609891
$ws.Range("O15").Value2 = 0.0000555027340683742
$ws.Range("P15").Value2 = 0.0000514023656946042
$ws.Range("Q15").Value2 = 0.0000213699378759098
$ws.Range("R15").Value2 = 0.00000162537124726021
$ws.Range("S15").Value2 = 0.00100551375796416
$ws.Range("T15").Value2 = 0.0000166415851566074
$ws.Range("U15").Value2 = 0.00000151455048040156
$ws.Range("V15").Value2 = 0.000000424812939624828
# row 16
$ws.Range("M16").Value2 = 0.0000592170416014638
$ws.Range("N16").Value2 = 0.0000398199634075985
$ws.Range("O16").Value2 = 0.0000833001927724528
$ws.Range("P16").Value2 = 0.0000557665998073628
$ws.Range("Q16").Value2 = 0.0000252954685581056
$ws.Range("R16").Value2 = 0.00000193504506020531
$ws.Range("S16").Value2 = 0.00136884621283006
$ws.Range("T16").Value2 = 0.0000117501531366684
$ws.Range("U16").Value2 = 0.00000121231738711658
$ws.Range("V16").Value2 = 0.000000233137959060881
# row 17
$ws.Range("M17").Value2 = 0.0000748848980359754
$ws.Range("N17").Value2 = 0.0000371794521487301
$ws.Range("O17").Value2 = 0.0000552570264250672
$ws.Range("P17").Value2 = 0.0000540112517812155
$ws.Range("Q17").Value2 = 0.0000266042096164778
$ws.Range("R17").Value2 = 0.00000193787166821378
$ws.Range("S17").Value2 = 0.00175003653580192
$ws.Range("T17").Value2 = 0.000022313208065433
$ws.Range("U17").Value2 = 0.0000019101877872393
$ws.Range("V17").Value2 = 0.000000193787166821378
# row 18
$ws.Range("M18").Value2 = 0.000101030549681197
$ws.Range("N18").Value2 = 0.0000503724956578947
$ws.Range("O18").Value2 = 0.000051971622504177
$ws.Range("P18").Value2 = 0.0000694858689158411
$ws.Range("Q18").Value2 = 0.0000270518958162768
$ws.Range("R18").Value2 = 0.00000182757353860842
$ws.Range("S18").Value2 = 0.00144513473843025
$ws.Range("T18").Value2 = 0.0000123361213856069
$ws.Range("U18").Value2 = 0.00000196083410913195
$ws.Range("V18").Value2 = 0.00000182757353860842
# row 19
$ws.Range("M19").Value2 = 0.000108590131906717
$ws.Range("N19").Value2 = 0.000033580277971605
$ws.Range("O19").Value2 = 0.0000779778697382997
$ws.Range("P19").Value2 = 0.0000721890124692722
$ws.Range("Q19").Value2 = 0.0000273008395780836
$ws.Range("R19").Value2 = 0.00000191326638552607
$ws.Range("S19").Value2 = 0.00118781954768077
$ws.Range("T19").Value2 = 0.0000109890171886626
$ws.Range("U19").Value2 = 0.0000014472143172569
$ws.Range("V19").Value2 = 0.000000735871686740797
# row 20
$ws.Range("M20").Value2 = 0.0000839875168832107
$ws.Range("N20").Value2 = 0.0000349348504380379
$ws.Range("O20").Value2 = 0.0000539309563712052
$ws.Range("P20").Value2 = 0.0000640472258030694
$ws.Range("Q20").Value2 = 0.0000244364108276365
$ws.Range("R20").Value2 = 0.00000200077328763537
$ws.Range("S20").Value2 = 0.00128368714910736
$ws.Range("T20").Value2 = 0.0000129713054715237
$ws.Range("U20").Value2 = 0.00000139379712172352
$ws.Range("V20").Value2 = 0.000000584495567174378
# row 21
$ws.Range("M21").Value2 = 0.000136812758068851
$ws.Range("N21").Value2 = 0.0000640638743311588
$ws.Range("O21").Value2 = 0.000111258788084243
$ws.Range("P21").Value2 = 0.0000738702997874368
$ws.Range("Q21").Value2 = 0.0000207104070463486
$ws.Range("R21").Value2 = 0.000001646334200689
$ws.Range("S21").Value2 = 0.00107751380439297
$ws.Range("T21").Value2 = 0.0000249097522539031
$ws.Range("U21").Value2 = 0.00000376986672041829
$ws.Range("V21").Value2 = 0.0000021473924356813
# row 22
$ws.Range("M22").Value2 = 0.0000759815245214799
$ws.Range("N22").Value2 = 0.0000405535565481625
$ws.Range("O22").Value2 = 0.0000264049269613654
$ws.Range("P22").Value2 = 0.0000424834847646078
$ws.Range("Q22").Value2 = 0.0000248760227898957
$ws.Range("R22").Value2 = 0.00000142864816022575
$ws.Range("S22").Value2 = 0.00118047693639285
$ws.Range("T22").Value2 = 0.00000199258822347275
$ws.Range("U22").Value2 = 0.00000187980021082335
$ws.Range("V22").Value2 = 0.000000726856081518363
# row 23
$ws.Range("M23").Value2 = 0.0000938803349907068
$ws.Range("N23").Value2 = 0.0000539102972951404
$ws.Range("O23").Value2 = 0.0000343133370657607
$ws.Range("P23").Value2 = 0.0000416416748209971
$ws.Range("Q23").Value2 = 0.0000250745568814606
$ws.Range("R23").Value2 = 0.00000176118911429307
$ws.Range("S23").Value2 = 0.00135227980781606
$ws.Range("T23").Value2 = 0.00000358207955449438
$ws.Range("U23").Value2 = 0.000017447712496683
$ws.Range("V23").Value2 = 0.00000155223447361423
# row 24
$ws.Range("M24").Value2 = 0.000100398353277247
$ws.Range("N24").Value2 = 0.0000479616011166941
$ws.Range("O24").Value2 = 0.0000342989491665379
$ws.Range("P24").Value2 = 0.0000414149137239109
$ws.Range("Q24").Value2 = 0.0000207944297621011
$ws.Range("R24").Value2 = 0.00000154969894805012
$ws.Range("S24").Value2 = 0.00107172751535866
$ws.Range("T24").Value2 = 0.00000972515156174311
$ws.Range("U24").Value2 = 0.00000551882584560707
$ws.Range("V24").Value2 = 0.000000806475983168941
# row 25
$ws.Range("M25").Value2 = 0.000119035696169698
$ws.Range("N25").Value2 = 0.0000514228414961797
$ws.Range("O25").Value2 = 0.000030309210715715
$ws.Range("P25").Value2 = 0.000048381783564837
$ws.Range("Q25").Value2 = 0.0000268481971652822
$ws.Range("R25").Value2 = 0.0000018101535305611
$ws.Range("S25").Value2 = 0.00127256689444094
$ws.Range("T25").Value2 = 0.00000577801006955102
$ws.Range("U25").Value2 = 0.0000041271500496793
$ws.Range("V25").Value2 = 0.000000897836151158304
# row 26
$ws.Range("M26").Value2 = 0.000135075234011024
$ws.Range("N26").Value2 = 0.0000539093233936396
$ws.Range("O26").Value2 = 0.0000666943077751362
$ws.Range("P26").Value2 = 0.0000610097788889007
$ws.Range("Q26").Value2 = 0.000012597582330302
$ws.Range("R26").Value2 = 0.00000122852455783111
$ws.Range("S26").Value2 = 0.00114279853063465
$ws.Range("T26").Value2 = 0.00000893283110694144
$ws.Range("U26").Value2 = 0.000000707963304512841
$ws.Range("V26").Value2 = 0.00000122852455783111
# row 27
$ws.Range("M27").Value2 = 0.0000568689864910942
$ws.Range("N27").Value2 = 0.0000344319489852141
$ws.Range("O27").Value2 = 0.000107526831171037
$ws.Range("P27").Value2 = 0.0000508287935371303
$ws.Range("Q27").Value2 = 0.0000154566258373841
$ws.Range("R27").Value2 = 0.00000121088773841258
$ws.Range("S27").Value2 = 0.00105901392453944
$ws.Range("T27").Value2 = 0.0000071940977399806
$ws.Range("U27").Value2 = 0.000000170948857187658
$ws.Range("V27").Value2 = 0.0000000712286904948575
# row 28
$ws.Range("M28").Value2 = 0.0000440742047396436
$ws.Range("N28").Value2 = 0.000141375607415511
$ws.Range("O28").Value2 = 0.000222646559506812
$ws.Range("P28").Value2 = 0.0000498583879402603
$ws.Range("Q28").Value2 = 0.0000156999258302451
$ws.Range("R28").Value2 = 0.00000137294897948702
$ws.Range("S28").Value2 = 0.00167049753331918
$ws.Range("T28").Value2 = 0.0000132337026633888
$ws.Range("U28").Value2 = 0.000000381374716524173
$ws.Range("V28").Value2 = 0.000000177974867711281
# row 29
$ws.Range("M29").Value2 = 0.0000502570598891735
$ws.Range("N29").Value2 = 0.0000334129128760898
$ws.Range("O29").Value2 = 0.00016290324922372
$ws.Range("P29").Value2 = 0.000035279385115505
$ws.Range("Q29").Value2 = 0.0000180221663773048
$ws.Range("R29").Value2 = 0.00000133100889204204
$ws.Range("S29").Value2 = 0.0013744426189904
$ws.Range("T29").Value2 = 0.00000989842245001377
$ws.Range("U29").Value2 = 0.000000290680102859755
$ws.Range("V29").Value2 = 0.000000305979055641848
# row 30
$ws.Range("M30").Value2 = 0.0000849782097250331
$ws.Range("N30").Value2 = 0.0000495211012383643
$ws.Range("O30").Value2 = 0.0000768800531695829
$ws.Range("P30").Value2 = 0.0000388361953946984
$ws.Range("Q30").Value2 = 0.0000176062079932161
$ws.Range("R30").Value2 = 0.00000147980702524053
$ws.Range("S30").Value2 = 0.00155850479412616
$ws.Range("T30").Value2 = 0.0000104984734625332
$ws.Range("U30").Value2 = 0.000000722425476889079
$ws.Range("V30").Value2 = 0.000000198084404953457
# row 31
$ws.Range("M31").Value2 = 0.000145469523226999
$ws.Range("N31").Value2 = 0.0000502455475653045
$ws.Range("O31").Value2 = 0.0000176448502208642
$ws.Range("P31").Value2 = 0.0000514514407068718
$ws.Range("Q31").Value2 = 0.0000113797502209972
$ws.Range("R31").Value2 = 0.0000011643106194443
$ws.Range("S31").Value2 = 0.00095463768205937
$ws.Range("T31").Value2 = 0.00000279988982294938
$ws.Range("U31").Value2 = 0.000000776207079629532
$ws.Range("V31").Value2 = 0.000000277216814153404
# row 32
$ws.Range("M32").Value2 = 0.000169194061332965
$ws.Range("N32").Value2 = 0.0000770028072047026
$ws.Range("O32").Value2 = 0.000010514071447947
$ws.Range("P32").Value2 = 0.0000263440993191399
$ws.Range("Q32").Value2 = 0.0000517323739611937
$ws.Range("R32").Value2 = 0.00000282819356507665
$ws.Range("S32").Value2 = 0.00125210414000532
$ws.Range("T32").Value2 = 0.00000302459589598475
$ws.Range("U32").Value2 = 0.0000177285837366379
$ws.Range("V32").Value2 = 0.000000471365594179441
# row 33
$ws.Range("M33").Value2 = 0.000124440582820496
$ws.Range("N33").Value2 = 0.000044718478437831
$ws.Range("O33").Value2 = 0.0000319663882474669
$ws.Range("P33").Value2 = 0.0000480142044498677
$ws.Range("Q33").Value2 = 0.00000799159706186672
$ws.Range("R33").Value2 = 0.000000818546329786888
$ws.Range("S33").Value2 = 0.00102902044000656
$ws.Range("T33").Value2 = 0.0000172325543113029
$ws.Range("U33").Value2 = 0.000000667761479562988
$ws.Range("V33").Value2 = 0.000000323110393336929
# row 34
$ws.Range("M34").Value2 = 0.0000631422668610423
$ws.Range("N34").Value2 = 0.0000508000107323966
$ws.Range("O34").Value2 = 0.0000282131941807645
$ws.Range("P34").Value2 = 0.0000241804148396523
$ws.Range("Q34").Value2 = 0.0000167978268522938
$ws.Range("R34").Value2 = 0.00000136594138973154
$ws.Range("S34").Value2 = 0.00121637080755593
$ws.Range("T34").Value2 = 0.00000725249833119363
$ws.Range("U34").Value2 = 0.0000036099879585762
$ws.Range("V34").Value2 = 0.000000292701726371043
# row 35
$ws.Range("M35").Value2 = 0.000173012494569825
$ws.Range("N35").Value2 = 0.0000154461718968352
$ws.Range("O35").Value2 = 0.0000640593477273725
$ws.Range("P35").Value2 = 0.0000345656125333159
$ws.Range("Q35").Value2 = 0.0000115116246275916
$ws.Range("R35").Value2 = 0.000000999006142581384
$ws.Range("S35").Value2 = 0.00144657626378311
$ws.Range("T35").Value2 = 0.0000262969155377961
$ws.Range("U35").Value2 = 0.000000568665035007865
$ws.Range("V35").Value2 = 0.00000371937671545685

# --- E19 / E30: confirmed_forage_sp flips from "no" to "yes" ---
$ws.Range("E19").Value2 = "yes"
$ws.Range("E30").Value2 = "yes"
